$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$rows = @(
    @{ Row = 3;  A = "Olivier";    B = "Olivier Chambord";    C = "Morgan Lewis And Bockius LLP"; D = "France";      G = "olivier.chambord@morganlewis.com" },
    @{ Row = 5;  A = "James";      B = "James P Bradley";     C = "Morgan Lewis And Bockius LLP"; D = "Singapore";   G = "james.bradley@morganlewis.com" },
    @{ Row = 6;  A = "Alexandra";  B = "Alexandra Rodina";    C = "Kennedys";                     D = "England";     G = "alexandra.rodina@kennedyslaw.com" },
    @{ Row = 10; A = "Tomoko";     B = "Tomoko Fuminaga";     C = "Morgan Lewis And Bockius LLP"; D = "Japan";       G = "tomoko.fuminaga@morganlewis.com" },
    @{ Row = 12; A = "Amanda";     B = "Amanda Beaumont";     C = "Kennedys";                     D = "England";     G = "amanda.beaumont@kennedyslaw.com" },
    @{ Row = 13; A = "Alberto";    B = "Alberto Bunge";       C = "Kennedys";                     D = "Argentina";   G = "alberto.bunge@kennedyslaw.com" },
    @{ Row = 15; A = "Andrea";     B = "Andrea Dougall";      C = "Morgan Lewis And Bockius LLP"; D = "the UAE";     G = "andrea.dougall@morganlewis.com" },
    @{ Row = 16; A = "Bingna";     B = "Bingna Guo";          C = "Morgan Lewis And Bockius LLP"; D = "China";       G = "bingna.guo@morganlewis.com" },
    @{ Row = 17; A = "Alfonso";    B = "Alfonso De Ramos";    C = "Kennedys";                     D = "Spain";       G = "alfonso.deramos@kennedyslaw.com" },
    @{ Row = 23; A = "Adam";       B = "Adam Longney";        C = "Kennedys";                     D = "England";     G = "adam.longney@kennedyslaw.com" },
    @{ Row = 34; A = "Nick";       B = "Nick Bolter";         C = "Morgan Lewis And Bockius LLP"; D = "Belgium";     G = "nick.bolter@morganlewis.com" },
    @{ Row = 36; A = "Alex";       B = "Alex Nurse";          C = "Kennedys";                     D = "England";     G = "alex.nurse@kennedyslaw.com" },
    @{ Row = 37; A = "Andy";       B = "Andy Purssell";       C = "Kennedys";                     D = "England";     G = "andrew.purssell@kennedyslaw.com" },
    @{ Row = 42; A = "Alistair";   B = "Alistair Darroch";    C = "Kennedys";                     D = "New Zealand"; G = "alistair.darroch@kennedyslaw.com" },
    @{ Row = 47; A = "Alberto";    B = "Alberto Torres";      C = "Kennedys";                     D = "Mexico";      G = "alberto.torres@kennedyslaw.com" },
    @{ Row = 52; A = "Alexandre";  B = "Alexandre Bailly";    C = "Morgan Lewis And Bockius LLP"; D = "France";      G = "alexandre.bailly@morganlewis.com" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
